# Applies the "gh-pages output generated at 456a3b4" update:
#   1. Bumps several "want to go" (F column) visitor counts that increased
#      since the last scrape, on sheets "展览" (Exhibition), "演出"
#      (Performance) and "全部类型" (All types).
#   2. On sheet "展览", refreshes the event list around 2024-06/07:
#        - a new event ("北京·EXA·全职高手ONLY·夏令营", 2024-06-29) is
#          inserted at row 46,
#        - the event that used to sit at row 46 (ICOS×CGF 03) moves down
#          to row 47 (its own want-count also ticked up from 9 to 11),
#        - the event that used to sit at row 47 (IDO 46th) moves down to
#          row 48 (want-count 223 -> 224),
#        - the event that used to sit at row 48 (万游引力 s8) drops out of
#          the list entirely,
#        - row 49 (梦次元 M30) is untouched,
#        - row 50 (IDO 47th) keeps its place, only its want-count changes
#          (38 -> 40).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibition)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# B46 / B48 store plain-text dates ("2024-06-29" / "2024-07-20"), not real
# Excel date values. Mark the cells as Text *before* writing so the
# date-looking string isn't auto-converted into a date serial number.
$ws1.Range("B46").NumberFormat = "@"
$ws1.Range("B48").NumberFormat = "@"

# Simple want-count refreshes (row content unchanged).
$ws1.Range("F3").Value  = 2930
$ws1.Range("F4").Value  = 2930
$ws1.Range("F11").Value = 2924
$ws1.Range("F14").Value = 7315
$ws1.Range("F16").Value = 64
$ws1.Range("F18").Value = 239
$ws1.Range("F19").Value = 116
$ws1.Range("F21").Value = 8876
$ws1.Range("F24").Value = 267
$ws1.Range("F29").Value = 101
$ws1.Range("F33").Value = 64
$ws1.Range("F34").Value = 105
$ws1.Range("F35").Value = 2615
$ws1.Range("F39").Value = 1478
$ws1.Range("F40").Value = 730
$ws1.Range("F41").Value = 3844
$ws1.Range("F43").Value = 196

# Row 46: brand-new event inserted.
$ws1.Range("B46").Value = "2024-06-29"
$ws1.Range("C46").Value = "北京·EXA·全职高手ONLY·夏令营"
$ws1.Range("D46").Value = "金盏路6号 蓝可可亲子乐园(金盏店)"
$ws1.Range("E46").Value = "2024.06.29 10:00-06.29 16:00"
$ws1.Range("F46").Value = 4
$ws1.Range("G46").Value = 68
$ws1.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=83977"
$ws1.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202404/JykdQ3eR1712654574985.jpeg"

# Row 47: previous row-46 event (ICOS x CGF 03), shifted down one row.
$ws1.Range("C47").Value = " 北京·ICOS国际动漫节×CGF中国游戏节03"
$ws1.Range("D47").Value = "石景山路68号 北京首钢会展中心"
$ws1.Range("E47").Value = "2024.07.20 09:00-07.21 17:00"
$ws1.Range("F47").Value = 11
$ws1.Range("G47").Value = 70
$ws1.Range("H47").Value = "https://show.bilibili.com/platform/detail.html?id=83931"
$ws1.Range("I47").Value = "//i1.hdslb.com/bfs/openplatform/202404/sgFsCjWK1712558620744.jpeg"

# Row 48: previous row-47 event (IDO 46th), shifted down one row
# (replaces the 万游引力 s8 event that used to be here, now dropped).
$ws1.Range("B48").Value = "2024-07-20"
$ws1.Range("C48").Value = "北京·IDO动漫游戏嘉年华46th"
$ws1.Range("D48").Value = "京沈路与天北路交汇处西北角 中国国际展览中心新馆"
$ws1.Range("E48").Value = "2024.07.20 09:30-07.21 17:00"
$ws1.Range("F48").Value = 224
$ws1.Range("G48").Value = 75
$ws1.Range("H48").Value = "https://show.bilibili.com/platform/detail.html?id=83716"
$ws1.Range("I48").Value = "//i0.hdslb.com/bfs/openplatform/202404/G4DiYbc51712040520493.jpeg"

# Row 50 (IDO 47th): unchanged event, want-count refreshed.
$ws1.Range("F50").Value = 40

# ---------------------------------------------------------------------
# Sheet "演出" (Performance)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 46
$ws2.Range("F5").Value = 259
$ws2.Range("F8").Value = 27

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 46
$ws4.Range("F4").Value  = 2930
$ws4.Range("F7").Value  = 259
$ws4.Range("F15").Value = 27
$ws4.Range("F19").Value = 7315
$ws4.Range("F23").Value = 239
$ws4.Range("F24").Value = 116
$ws4.Range("F25").Value = 8876
$ws4.Range("F31").Value = 101
$ws4.Range("F35").Value = 64
$ws4.Range("F37").Value = 105
$ws4.Range("F38").Value = 2615
$ws4.Range("F42").Value = 730
$ws4.Range("F43").Value = 3844
$ws4.Range("F44").Value = 196
$ws4.Range("F48").Value = 224
$ws4.Range("F50").Value = 40
